$d = $word.ActiveDocument

function Escape-Xml($s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    $s = $s -replace '"', "&quot;"
    return $s
}

# --- Locate the two paragraphs to rewrite -------------------------------
# Para 1 currently: Heading1-styled "Day After Day - May 1935"
# Para 2 currently: bold "By Dorothy Day"
# (Title/Authors styles already exist in styles.xml for this template.)
$titlePara = $d.Paragraphs.Item(1)
$authorsPara = $d.Paragraphs.Item(2)

if ($titlePara.Range.Text.TrimEnd([char]13) -ne "Day After Day - May 1935") {
    throw "Unexpected content in paragraph 1: $($titlePara.Range.Text)"
}
if ($authorsPara.Range.Text.TrimEnd([char]13) -ne "By Dorothy Day") {
    throw "Unexpected content in paragraph 2: $($authorsPara.Range.Text)"
}

# --- Rebuild the title as a pandoc-style title block ---------------------
# "Day After Day - May 1935" -> Title-styled paragraph, one run per
# word/separator (mirrors how a pandoc title-block importer emits runs).
$titleWords = @("Day", " ", "After", " ", "Day", " ", "-", " ", "May", " ", "1935")
$titleRuns = ($titleWords | ForEach-Object {
    '<w:r><w:t xml:space="preserve">' + (Escape-Xml $_) + '</w:t></w:r>'
}) -join ""
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>'
[void]$titlePara.Range.InsertXML($titleXml)

# --- Rebuild the byline as the pandoc "author" block ----------------------
# "By Dorothy Day" -> Authors-styled paragraph "Dorothy Day" (drop "By "),
# again one run per word/separator.
$authorWords = @("Dorothy", " ", "Day")
$authorRuns = ($authorWords | ForEach-Object {
    '<w:r><w:t xml:space="preserve">' + (Escape-Xml $_) + '</w:t></w:r>'
}) -join ""
$authorsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'
[void]$authorsPara.Range.InsertXML($authorsXml)

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
